$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update existing values
$ws.Range("D2").Value = 2772
$ws.Range("E2").Value = 141
$ws.Range("F2").Value = 141
$ws.Range("G2").Value = 129
$ws.Range("H2").Value = 101
$ws.Range("I2").Value = 99
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 3130
$ws.Range("L2").Value = 1148
$ws.Range("M2").Value = 1982
$ws.Range("N2").Value = 1855
$ws.Range("O2").Value = 127
$ws.Range("P2").Value = 202
$ws.Range("Q2").Value = 184
$ws.Range("R2").Value = -122
$ws.Range("S2").Value = -27
$ws.Range("T2").Value = 192
$ws.Range("U2").Value = -9
$ws.Range("V2").Value = 724
$ws.Range("W2").Value = 5.08
$ws.Range("X2").Value = 3.66
$ws.Range("Y2").Value = 5.48
$ws.Range("Z2").Value = 3.22
$ws.Range("AA2").Value = 57.91
$ws.Range("AB2").Value = 832.48
$ws.Range("AC2").Value = 246
$ws.Range("AD2").Value = 17.34
$ws.Range("AE2").Value = 4754
$ws.Range("AF2").Value = 0.9
$ws.Range("AG2").Value = 50
$ws.Range("AH2").Value = 1.17
$ws.Range("AI2").Value = 19.66
$ws.Range("AJ2").Value = 40396365

# Row 3: update existing values
$ws.Range("D3").Value = 2632
$ws.Range("E3").Value = 67
$ws.Range("F3").Value = 67
$ws.Range("G3").Value = 58
$ws.Range("H3").Value = 37
$ws.Range("I3").Value = 19
$ws.Range("J3").Value = 18
$ws.Range("K3").Value = 3267
$ws.Range("L3").Value = 1256
$ws.Range("M3").Value = 2011
$ws.Range("N3").Value = 1867
$ws.Range("O3").Value = 145
$ws.Range("P3").Value = 202
$ws.Range("Q3").Value = 54
$ws.Range("R3").Value = -101
$ws.Range("S3").Value = 55
$ws.Range("T3").Value = 105
$ws.Range("U3").Value = -50
$ws.Range("V3").Value = 800
$ws.Range("W3").Value = 2.55
$ws.Range("X3").Value = 1.4
$ws.Range("Y3").Value = 1.01
$ws.Range("Z3").Value = 1.15
$ws.Range("AA3").Value = 62.43
$ws.Range("AB3").Value = 830.92
$ws.Range("AC3").Value = 47
$ws.Range("AD3").Value = 78.51
$ws.Range("AE3").Value = 4782
$ws.Range("AF3").Value = 0.77
$ws.Range("AG3").Value = 50
$ws.Range("AH3").Value = 1.36
$ws.Range("AI3").Value = 103.36
$ws.Range("AJ3").Value = 40396365

# Row 4: update existing values
$ws.Range("D4").Value = 2749
$ws.Range("E4").Value = 106
$ws.Range("F4").Value = 106
$ws.Range("G4").Value = 103
$ws.Range("H4").Value = 68
$ws.Range("I4").Value = 51
$ws.Range("J4").Value = 17
$ws.Range("K4").Value = 3443
$ws.Range("L4").Value = 1360
$ws.Range("M4").Value = 2083
$ws.Range("N4").Value = 1924
$ws.Range("O4").Value = 159
$ws.Range("P4").Value = 202
$ws.Range("Q4").Value = 190
$ws.Range("R4").Value = -70
$ws.Range("S4").Value = -73
$ws.Range("T4").Value = 86
$ws.Range("U4").Value = 104
$ws.Range("V4").Value = 727
$ws.Range("W4").Value = 3.85
$ws.Range("X4").Value = 2.49
$ws.Range("Y4").Value = 2.71
$ws.Range("Z4").Value = 2.04
$ws.Range("AA4").Value = 65.3
$ws.Range("AB4").Value = 852.14
$ws.Range("AC4").Value = 127
$ws.Range("AD4").Value = 39.32
$ws.Range("AE4").Value = 4867
$ws.Range("AF4").Value = 1.03
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 1
$ws.Range("AI4").Value = 38.56
$ws.Range("AJ4").Value = 40396365

# Row 5: update existing values
$ws.Range("D5").Value = 3000
$ws.Range("E5").Value = 74
$ws.Range("F5").Value = 74
$ws.Range("G5").Value = 63
$ws.Range("H5").Value = 48
$ws.Range("I5").Value = 36
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 3580
$ws.Range("L5").Value = 1470
$ws.Range("M5").Value = 2109
$ws.Range("N5").Value = 1944
$ws.Range("O5").Value = 166
$ws.Range("P5").Value = 202
$ws.Range("Q5").Value = 53
$ws.Range("R5").Value = -254
$ws.Range("S5").Value = 219
$ws.Range("T5").Value = 283
$ws.Range("U5").Value = -230
$ws.Range("V5").Value = 957
$ws.Range("W5").Value = 2.46
$ws.Range("X5").Value = 1.59
$ws.Range("Y5").Value = 1.85
$ws.Range("Z5").Value = 1.36
$ws.Range("AA5").Value = 69.7
$ws.Range("AB5").Value = 861.15
$ws.Range("AC5").Value = 89
$ws.Range("AD5").Value = 65.51
$ws.Range("AE5").Value = 4898
$ws.Range("AF5").Value = 1.18
$ws.Range("AG5").Value = 50
$ws.Range("AH5").Value = 0.86
$ws.Range("AI5").Value = 55.48
$ws.Range("AJ5").Value = 40396365

# Row 6: update existing values
$ws.Range("D6").Value = 3311
$ws.Range("E6").Value = 197
$ws.Range("F6").Value = 197
$ws.Range("G6").Value = 177
$ws.Range("H6").Value = 138
$ws.Range("I6").Value = 132
$ws.Range("K6").Value = 4130
$ws.Range("L6").Value = 1892
$ws.Range("M6").Value = 2238
$ws.Range("N6").Value = 2068
$ws.Range("P6").Value = 202
$ws.Range("Q6").Value = 19
$ws.Range("R6").Value = -317
$ws.Range("S6").Value = 315
$ws.Range("T6").Value = 182
$ws.Range("U6").Value = -163
$ws.Range("V6").Value = 1274
$ws.Range("W6").Value = 5.95
$ws.Range("X6").Value = 4.17
$ws.Range("Y6").Value = 6.59
$ws.Range("Z6").Value = 3.58
$ws.Range("AA6").Value = 84.57
$ws.Range("AB6").Value = 921.41
$ws.Range("AC6").Value = 327
$ws.Range("AD6").Value = 17.73
$ws.Range("AE6").Value = 5151
$ws.Range("AF6").Value = 1.13
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 0.86
$ws.Range("AI6").Value = 15.19
$ws.Range("AJ6").Value = 40396365

# Rows 7-9: clear out all data columns (D through AJ), keep A/B/C labels
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
